# Apply updated odds values to Sheet1, as captured in the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("M3").Value = 1.13
$ws.Range("O3").Value = 1.62
$ws.Range("P3").Value = 2.2
$ws.Range("W3").Value = 6.5
$ws.Range("X3").Value = 1.11

# Row 4
$ws.Range("M4").Value = 1.14
$ws.Range("O4").Value = 1.67
$ws.Range("X4").Value = 1.1
$ws.Range("Z4").Value = 2.08

# Row 5
$ws.Range("G5").Value = 1.57
$ws.Range("H5").Value = 3.6
$ws.Range("I5").Value = 6.5
$ws.Range("J5").Value = 2.3
$ws.Range("K5").Value = 1.91
$ws.Range("L5").Value = 9
$ws.Range("M5").Value = 1.14
$ws.Range("N5").Value = 4.9
$ws.Range("O5").Value = 1.62
$ws.Range("P5").Value = 2.15
$ws.Range("S5").Value = 3.1
$ws.Range("T5").Value = 1.36
$ws.Range("W5").Value = 5.3
$ws.Range("X5").Value = 1.12
$ws.Range("AA5").Value = 2.87
$ws.Range("AB5").Value = 1.37
$ws.Range("AC5").Value = 4
$ws.Range("AD5").Value = 5.5
$ws.Range("AF5").Value = 11
$ws.Range("AJ5").Value = 8.5
$ws.Range("AK5").Value = 41
$ws.Range("AL5").Value = 201
$ws.Range("AN5").Value = 9.5
$ws.Range("AO5").Value = 29
$ws.Range("AP5").Value = 23
$ws.Range("AQ5").Value = 101
$ws.Range("AR5").Value = 81
$ws.Range("AS5").Value = 101

# Row 7
$ws.Range("I7").Value = 5.5
$ws.Range("J7").Value = 2.22
$ws.Range("N7").Value = 6.5
$ws.Range("O7").Value = 1.36
$ws.Range("P7").Value = 2.92
$ws.Range("S7").Value = 2.05
$ws.Range("T7").Value = 1.7
$ws.Range("W7").Value = 3.5
$ws.Range("X7").Value = 1.26
$ws.Range("AA7").Value = 1.98
$ws.Range("AB7").Value = 1.75
$ws.Range("AC7").Value = 5.7
$ws.Range("AD7").Value = 7
$ws.Range("AF7").Value = 12.5
$ws.Range("AG7").Value = 14.5
$ws.Range("AH7").Value = 30
$ws.Range("AI7").Value = 6.5
$ws.Range("AK7").Value = 17
$ws.Range("AL7").Value = 90
$ws.Range("AM7").Value = 800
$ws.Range("AN7").Value = 13.5
$ws.Range("AP7").Value = 17
$ws.Range("AR7").Value = 60
$ws.Range("AS7").Value = 60

# Row 9
$ws.Range("G9").Value = 1.4
$ws.Range("H9").Value = 3.85
$ws.Range("I9").Value = 9.25
$ws.Range("J9").Value = 1.88
$ws.Range("K9").Value = 2.18
$ws.Range("L9").Value = 8.25
$ws.Range("M9").Value = 1.08
$ws.Range("N9").Value = 6.5
$ws.Range("O9").Value = 1.36
$ws.Range("P9").Value = 2.9
$ws.Range("S9").Value = 2.05
$ws.Range("T9").Value = 1.7
$ws.Range("W9").Value = 3.5
$ws.Range("X9").Value = 1.26
$ws.Range("Y9").Value = 1.42
$ws.Range("Z9").Value = 2.67
$ws.Range("AA9").Value = 2.32
$ws.Range("AB9").Value = 1.55
$ws.Range("AC9").Value = 5.1
$ws.Range("AF9").Value = 8.75
$ws.Range("AG9").Value = 13.5
$ws.Range("AI9").Value = 6.5
$ws.Range("AJ9").Value = 7.9
$ws.Range("AK9").Value = 25
$ws.Range("AN9").Value = 18
$ws.Range("AP9").Value = 29
$ws.Range("AR9").Value = 150
$ws.Range("AS9").Value = 120

# Row 10
$ws.Range("G10").Value = 1.91
$ws.Range("H10").Value = 3.6
$ws.Range("I10").Value = 3.3
$ws.Range("J10").Value = 2.5
$ws.Range("L10").Value = 3.75
$ws.Range("M10").Value = 1.03
$ws.Range("O10").Value = 1.17
$ws.Range("W10").Value = 2.38
$ws.Range("X10").Value = 1.53
$ws.Range("AA10").Value = 1.5
$ws.Range("AB10").Value = 2.5
$ws.Range("AC10").Value = 11
$ws.Range("AD10").Value = 12
$ws.Range("AE10").Value = 9
$ws.Range("AF10").Value = 19
$ws.Range("AG10").Value = 15
$ws.Range("AJ10").Value = 7.5
$ws.Range("AK10").Value = 12
$ws.Range("AL10").Value = 34
$ws.Range("AM10").Value = 101
$ws.Range("AO10").Value = 21
$ws.Range("AP10").Value = 12
$ws.Range("AR10").Value = 23
$ws.Range("AS10").Value = 26

# Row 12
$ws.Range("K12").Value = 2.75
$ws.Range("M12").Value = 1.02
$ws.Range("N12").Value = 11
$ws.Range("O12").Value = 1.14
$ws.Range("X12").Value = 1.62
$ws.Range("AA12").Value = 2
$ws.Range("AB12").Value = 1.73
$ws.Range("AC12").Value = 9
$ws.Range("AE12").Value = 10
$ws.Range("AM12").Value = 1000

# Row 13
$ws.Range("M13").Value = 1.03
$ws.Range("O13").Value = 1.17
$ws.Range("W13").Value = 2.38
$ws.Range("X13").Value = 1.53

# Row 14
$ws.Range("M14").Value = 1.04
$ws.Range("O14").Value = 1.25
$ws.Range("X14").Value = 1.36

# Row 15
$ws.Range("G15").Value = 1.7
$ws.Range("H15").Value = 4
$ws.Range("J15").Value = 2.25
$ws.Range("K15").Value = 2.38
$ws.Range("L15").Value = 4.75
$ws.Range("M15").Value = 1.03
$ws.Range("N15").Value = 15
$ws.Range("O15").Value = 1.2
$ws.Range("P15").Value = 4.33
$ws.Range("S15").Value = 1.65
$ws.Range("T15").Value = 2.2
$ws.Range("U15").Value = 2.05
$ws.Range("V15").Value = 1.8
$ws.Range("W15").Value = 2.63
$ws.Range("X15").Value = 1.44
$ws.Range("Y15").Value = 1.3
$ws.Range("Z15").Value = 3.4
$ws.Range("AA15").Value = 1.62
$ws.Range("AB15").Value = 2.2
$ws.Range("AC15").Value = 9
$ws.Range("AD15").Value = 9
$ws.Range("AF15").Value = 13
$ws.Range("AG15").Value = 13
$ws.Range("AH15").Value = 21
$ws.Range("AI15").Value = 15
$ws.Range("AJ15").Value = 8
$ws.Range("AL15").Value = 41
$ws.Range("AM15").Value = 151
$ws.Range("AN15").Value = 15
$ws.Range("AO15").Value = 26
$ws.Range("AR15").Value = 34
$ws.Range("AS15").Value = 34

# Row 16
$ws.Range("G16").Value = 3.8
$ws.Range("I16").Value = 1.95
$ws.Range("M16").Value = 1.05
$ws.Range("N16").Value = 11
$ws.Range("O16").Value = 1.25
$ws.Range("S16").Value = 1.85
$ws.Range("T16").Value = 2
$ws.Range("X16").Value = 1.36

# Row 17
$ws.Range("G17").Value = 6
$ws.Range("H17").Value = 4.25
$ws.Range("I17").Value = 1.47
$ws.Range("J17").Value = 5.6
$ws.Range("K17").Value = 2.35
$ws.Range("L17").Value = 1.98
$ws.Range("O17").Value = 1.21
$ws.Range("P17").Value = 3.9
$ws.Range("T17").Value = 2.12
$ws.Range("W17").Value = 2.55
$ws.Range("AA17").Value = 1.78
$ws.Range("AB17").Value = 1.93
$ws.Range("AF17").Value = 120
$ws.Range("AH17").Value = 55
$ws.Range("AQ17").Value = 10.25
$ws.Range("AR17").Value = 11.25
$ws.Range("AS17").Value = 24

# Row 18
$ws.Range("G18").Value = 2.05
$ws.Range("H18").Value = 3.4
$ws.Range("I18").Value = 3.25
$ws.Range("J18").Value = 2.62
$ws.Range("L18").Value = 3.8
$ws.Range("Y18").Value = 1.39
$ws.Range("Z18").Value = 2.55
$ws.Range("AA18").Value = 1.82
$ws.Range("AB18").Value = 1.78
$ws.Range("AC18").Value = 6.8
$ws.Range("AD18").Value = 9.25
$ws.Range("AE18").Value = 8.75
$ws.Range("AF18").Value = 18
$ws.Range("AG18").Value = 17.5
$ws.Range("AJ18").Value = 6.6
$ws.Range("AK18").Value = 16
$ws.Range("AN18").Value = 9.25
$ws.Range("AO18").Value = 16.5
$ws.Range("AP18").Value = 11.75
$ws.Range("AQ18").Value = 45
$ws.Range("AR18").Value = 30
